$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '72.452.23'
$r.ClearFormats()
$ws.Range("E2").Value = '  +6.08%  '
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '4.062.27'
$r.ClearFormats()
$ws.Range("E3").Value = '  +6.39%  '
$ws.Range("E4").Value = '  +0.10%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '525.78'
$r.ClearFormats()
$ws.Range("E5").Value = '  +2.25%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '148.88'
$r.ClearFormats()
$ws.Range("E6").Value = '  +7.48%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.710'
$r.ClearFormats()
$ws.Range("E7").Value = '  +18.52%  '
$ws.Range("E8").Value = '  +0.09%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.763'
$r.ClearFormats()
$ws.Range("E9").Value = '  +8.92%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.177'
$r.ClearFormats()
$ws.Range("E10").Value = '  +7.78%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.0000336'
$r.ClearFormats()
$ws.Range("E11").Value = '  +6.71%  '
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '48.72'
$r.ClearFormats()
$ws.Range("E12").Value = '  +18.75%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '10.94'
$r.ClearFormats()
$ws.Range("E13").Value = '  +8.16%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '4.706.18'
$r.ClearFormats()
$ws.Range("E14").Value = '  +6.56%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '4.083.92'
$r.ClearFormats()
$ws.Range("E15").Value = '  +7.01%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '14.39'
$r.ClearFormats()
$ws.Range("E16").Value = '  +2.36%  '
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '21.08'
$r.ClearFormats()
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("E18").Value = '  +3.54%  '
$ws.Range("E19").Value = '  +0.43%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '72.375.68'
$r.ClearFormats()
$ws.Range("E20").Value = '  +6.17%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '440.21'
$r.ClearFormats()
$ws.Range("E21").Value = '  +6.66%  '
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '102.14'
$r.ClearFormats()
$ws.Range("E22").Value = '  +19.18%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '3.64'
$r.ClearFormats()
$ws.Range("E23").Value = '  +8.40%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '15.04'
$r.ClearFormats()
$ws.Range("E24").Value = '  +8.73%  '
$ws.Range("E25").Value = '  +6.81%  '
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '11.42'
$r.ClearFormats()
$ws.Range("E26").Value = '  -0.27%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '11.09'
$r.ClearFormats()
$ws.Range("E27").Value = '  +7.32%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '37.65'
$r.ClearFormats()
$ws.Range("E28").Value = '  +7.76%  '
$ws.Range("E29").Value = '  +3.38%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '3.44'
$r.ClearFormats()
$ws.Range("E30").Value = '  +23.16%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '13.70'
$r.ClearFormats()
$ws.Range("E31").Value = '  +5.40%  '
$ws.Range("E32").Value = '  +8.07%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '676.40'
$r.ClearFormats()
$ws.Range("E33").Value = '  +0.56%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '6.75'
$r.ClearFormats()
$ws.Range("E34").Value = '  +10.96%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '67.22'
$r.ClearFormats()
$ws.Range("E35").Value = '  +3.74%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '42.71'
$r.ClearFormats()
$ws.Range("E36").Value = '  +8.58%  '
$ws.Range("E37").Value = '  +7.66%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.434'
$r.ClearFormats()
$ws.Range("E38").Value = '  +0.01%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.158'
$r.ClearFormats()
$ws.Range("E39").Value = '  +6.48%  '
$ws.Range("E40").Value = '  +6.32%  '
$ws.Range("E41").Value = '  +8.79%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("E43").Value = '  -0.07%  '
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '3.15'
$r.ClearFormats()
$ws.Range("E44").Value = '  +2.64%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.155'
$r.ClearFormats()
$ws.Range("E45").Value = '  +12.91%  '
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("E47").Value = '  +3.51%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '9.50'
$r.ClearFormats()
$ws.Range("E48").Value = '  +13.06%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '3.11'
$r.ClearFormats()
$ws.Range("E49").Value = '  +6.89%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.000279'
$r.ClearFormats()
$ws.Range("E50").Value = '  +9.63%  '
$ws.Range("E51").Value = '  +5.19%  '
